# Feature: add arrows (arrow_n). Fixed bugs, removed unnecessary code.
#
# On the "meta" sheet, insert a new "style" / "default" key-value row right
# after the existing "bar_stack_index" / "1" row (i.e. as the new row 6),
# pushing the old (empty) trailing row down to row 7.

$wb = $excel.ActiveWorkbook

# First sheet in the workbook is the "meta" sheet.
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at row 6 (shifts the old row 6 -> row 7, and
# copies formatting from the row above, keeping column A bold/colored).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row with the style metadata key/value pair.
$ws.Range("A6").Value = "style"
$ws.Range("B6").Value = "default"
